# Generate Report for Handback
# This script updates the localization-status workbook to reflect that the
# zh-cn and de-de handoffs have now been handed back (i.e. the target files
# are in sync with en-US): the Status column is updated, the "Latest Target
# File" / "Latest Handback File" columns are populated (mirroring the source
# file name / handoff file name, as hyperlinks), and a "Latest Handback
# DateTime" value is recorded.

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# Overview sheet: the Status values shown here mirror the per-language
# sheets, so they need to be refreshed too.
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B2").Value = $newStatus
$overview.Range("C2").Value = $newStatus
$overview.Range("B3").Value = $newStatus
$overview.Range("C3").Value = $newStatus

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

$zh.Range("B2").Value = $newStatus
$zh.Range("B3").Value = $newStatus

$zh.Hyperlinks.Add($zh.Range("E2"), "https://github.com/OpenLocalizationTest/oltest/blob/334cf836b6b08200cdc3b286074aecf042c62c42/e2e/6685b6bf-1f52-4832-87df-291ee63b83d0.md", $null, $null, "6685b6bf-1f52-4832-87df-291ee63b83d0.md")
$zh.Hyperlinks.Add($zh.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/1d5fd25ece19f2c5bc5e6894c73b60c70c15484b/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/yuwzho/6685b6bf-1f52-4832-87df-291ee63b83d0.86fa7517248cbe8736fda64f533993182afad7b5.zh-cn.xlf", $null, $null, "6685b6bf-1f52-4832-87df-291ee63b83d0.86fa7517248cbe8736fda64f533993182afad7b5.zh-cn.xlf")
$zh.Range("G2").Value = "2016-01-25 08:40:51"

$zh.Hyperlinks.Add($zh.Range("E3"), "https://github.com/OpenLocalizationTest/oltest/blob/334cf836b6b08200cdc3b286074aecf042c62c42/e2e/bc478b3f-5523-40df-8c0e-efab388c8c5f.md", $null, $null, "bc478b3f-5523-40df-8c0e-efab388c8c5f.md")
$zh.Hyperlinks.Add($zh.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/1d5fd25ece19f2c5bc5e6894c73b60c70c15484b/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/yuwzho/bc478b3f-5523-40df-8c0e-efab388c8c5f.64ea431d83372bb592cc1fde2022869eebde10aa.zh-cn.xlf", $null, $null, "bc478b3f-5523-40df-8c0e-efab388c8c5f.64ea431d83372bb592cc1fde2022869eebde10aa.zh-cn.xlf")
$zh.Range("G3").Value = "2016-01-25 08:40:51"

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

$de.Range("B2").Value = $newStatus
$de.Range("B3").Value = $newStatus

$de.Hyperlinks.Add($de.Range("E2"), "https://github.com/OpenLocalizationTest/oltest/blob/334cf836b6b08200cdc3b286074aecf042c62c42/e2e/6685b6bf-1f52-4832-87df-291ee63b83d0.md", $null, $null, "6685b6bf-1f52-4832-87df-291ee63b83d0.md")
$de.Hyperlinks.Add($de.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/e6a660e1cec4260877b97354e96ca65993244249/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/yuwzho/6685b6bf-1f52-4832-87df-291ee63b83d0.86fa7517248cbe8736fda64f533993182afad7b5.de-de.xlf", $null, $null, "6685b6bf-1f52-4832-87df-291ee63b83d0.86fa7517248cbe8736fda64f533993182afad7b5.de-de.xlf")
$de.Range("G2").Value = "2016-01-25 08:41:14"

$de.Hyperlinks.Add($de.Range("E3"), "https://github.com/OpenLocalizationTest/oltest/blob/334cf836b6b08200cdc3b286074aecf042c62c42/e2e/bc478b3f-5523-40df-8c0e-efab388c8c5f.md", $null, $null, "bc478b3f-5523-40df-8c0e-efab388c8c5f.md")
$de.Hyperlinks.Add($de.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/e6a660e1cec4260877b97354e96ca65993244249/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/yuwzho/bc478b3f-5523-40df-8c0e-efab388c8c5f.64ea431d83372bb592cc1fde2022869eebde10aa.de-de.xlf", $null, $null, "bc478b3f-5523-40df-8c0e-efab388c8c5f.64ea431d83372bb592cc1fde2022869eebde10aa.de-de.xlf")
$de.Range("G3").Value = "2016-01-25 08:41:14"
